# Fixed update to excel issue
# Updates a couple of forecast numbers on "Forecast Comparison" and the
# rolled-up totals on "Summary" so they stay in sync.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Forecast Comparison: MyForecast (column D) corrections
$wsForecast.Range("D8").Value = 18
$wsForecast.Range("D16").Value = 13

# Summary: roll-up totals are stored as text, keep them that way
$wsSummary.Range("B9").Value = "'236"
$wsSummary.Range("B10").Value = "'139"
